# "Tried new logic for material demand"
#
# The reuse_scenario values on Sheet1 (column B) are reshuffled:
#   "LFP reused"       -> "No reuse"
#   "Direct recycling" -> "LFP reused"
# ("All reused" is untouched.)
#
# As a consequence of the relabeling, the LFP-chemistry "value" flags
# (column E) for the two affected scenario blocks are flipped so the
# data still lines up with the new scenario names:
#   - the block that used to be "LFP reused" / LFP (now "No reuse" / LFP)
#     goes from 1 -> 0 for all years (1950-2050)
#   - the block that used to be "Direct recycling" / LFP (now "LFP reused"
#     / LFP) goes from 0 -> 1 for all years (1950-2050)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column layout: A=index, B=reuse_scenario, C=battery_chemistry, D=year, E=value
# Data rows: 2-1718 = "LFP reused", 1719-3435 = "Direct recycling", 3436-5152 = "All reused"
# (17 chemistries x 101 years = 1717 rows per scenario block)

# Step 1: rename the "LFP reused" scenario block to "No reuse"
for ($r = 2; $r -le 1718; $r++) {
    $ws.Cells.Item($r, 2).Value = "No reuse"
}

# Step 2: rename the "Direct recycling" scenario block to "LFP reused"
for ($r = 1719; $r -le 3435; $r++) {
    $ws.Cells.Item($r, 2).Value = "LFP reused"
}

# Step 3: within the (now "No reuse") block, the LFP chemistry rows
# (years 1950-2050) flip value 1 -> 0
for ($r = 204; $r -le 304; $r++) {
    $ws.Cells.Item($r, 5).Value = 0
}

# Step 4: within the (now "LFP reused") block, the LFP chemistry rows
# (years 1950-2050) flip value 0 -> 1
for ($r = 1921; $r -le 2021; $r++) {
    $ws.Cells.Item($r, 5).Value = 1
}
